# Applies the cryptos-list refresh described by the commit diff.
# Each cell value is forced to literal text via a leading apostrophe
# (mirrors the workbook's original inlineStr/text storage) so Excel
# does not auto-coerce number-like strings (e.g. "30.551.71",
# "0.9990", "9.330") into numeric values and drop formatting like
# trailing zeros or multi-dot thousands separators.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.551.71"
$ws.Range("E2").Value = "'  +1.49%  "
$ws.Range("D3").Value = "'1.884.67"
$ws.Range("E3").Value = "'  +1.52%  "
$ws.Range("E4").Value = "'  -0.17%  "
$ws.Range("D5").Value = "'247.21"
$ws.Range("E5").Value = "'  +5.83%  "
$ws.Range("E6").Value = "'  -0.17%  "
$ws.Range("D7").Value = "'0.4747"
$ws.Range("E7").Value = "'  +1.23%  "
$ws.Range("D8").Value = "'0.2909"
$ws.Range("E8").Value = "'  +2.93%  "
$ws.Range("D9").Value = "'0.06533"
$ws.Range("E9").Value = "'  +1.37%  "
$ws.Range("D10").Value = "'22.11"
$ws.Range("E10").Value = "'  +5.43%  "
$ws.Range("D11").Value = "'0.07735"
$ws.Range("E11").Value = "'  -0.04%  "
$ws.Range("B12").Value = "'Polygon"
$ws.Range("C12").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.7447"
$ws.Range("E12").Value = "'  +9.82%  "
$ws.Range("B13").Value = "'Litecoin"
$ws.Range("C13").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").Value = "'97.44"
$ws.Range("E13").Value = "'  +4.34%  "
$ws.Range("D14").Value = "'1.878.00"
$ws.Range("E14").Value = "'  +1.10%  "
$ws.Range("D15").Value = "'5.167"
$ws.Range("E15").Value = "'  +2.46%  "
$ws.Range("D16").Value = "'275.28"
$ws.Range("E16").Value = "'  +3.68%  "
$ws.Range("D17").Value = "'30.539.09"
$ws.Range("E17").Value = "'  +1.49%  "
$ws.Range("D18").Value = "'13.64"
$ws.Range("E18").Value = "'  +2.55%  "
$ws.Range("D19").Value = "'0.000007584"
$ws.Range("E19").Value = "'  +0.24%  "
$ws.Range("D20").Value = "'0.9994"
$ws.Range("D21").Value = "'2.122.21"
$ws.Range("E21").Value = "'  +0.23%  "
$ws.Range("B22").Value = "'Uniswap"
$ws.Range("C22").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.285"
$ws.Range("E22").Value = "'  +2.86%  "
$ws.Range("B23").Value = "'BinanceUSD"
$ws.Range("C23").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "'0.9996"
$ws.Range("E23").Value = "'  -0.14%  "
$ws.Range("D24").Value = "'6.209"
$ws.Range("E24").Value = "'  +1.94%  "
$ws.Range("E25").Value = "'  +0.44%  "
$ws.Range("D26").Value = "'163.28"
$ws.Range("E26").Value = "'  -1.21%  "
$ws.Range("D27").Value = "'18.96"
$ws.Range("E27").Value = "'  +2.72%  "
$ws.Range("D28").Value = "'1.951"
$ws.Range("E28").Value = "'  +3.85%  "
$ws.Range("D29").Value = "'1.370"
$ws.Range("E29").Value = "'  +0.72%  "
$ws.Range("D30").Value = "'0.09985"
$ws.Range("E30").Value = "'  +1.54%  "
$ws.Range("D31").Value = "'1.520"
$ws.Range("E31").Value = "'  +4.96%  "
$ws.Range("D32").Value = "'4.331"
$ws.Range("E32").Value = "'  +3.42%  "
$ws.Range("D33").Value = "'4.088"
$ws.Range("E33").Value = "'  +3.00%  "
$ws.Range("D34").Value = "'0.04808"
$ws.Range("E34").Value = "'  +3.62%  "
$ws.Range("E35").Value = "'  +1.80%  "
$ws.Range("D36").Value = "'0.7035"
$ws.Range("E36").Value = "'  +2.33%  "
$ws.Range("D37").Value = "'2.717"
$ws.Range("E37").Value = "'  +0.07%  "
$ws.Range("D38").Value = "'0.01875"
$ws.Range("E38").Value = "'  +2.46%  "
$ws.Range("D39").Value = "'2.737"
$ws.Range("E39").Value = "'  +0.68%  "
$ws.Range("D40").Value = "'6.354"
$ws.Range("E40").Value = "'  +1.18%  "
$ws.Range("D41").Value = "'1.966"
$ws.Range("E41").Value = "'  +5.17%  "
$ws.Range("D42").Value = "'71.34"
$ws.Range("E42").Value = "'  +1.07%  "
$ws.Range("D43").Value = "'0.4239"
$ws.Range("E43").Value = "'  +4.99%  "
$ws.Range("D44").Value = "'0.8408"
$ws.Range("E44").Value = "'  +1.11%  "
$ws.Range("D45").Value = "'0.9990"
$ws.Range("E45").Value = "'  -0.15%  "
$ws.Range("D46").Value = "'102.89"
$ws.Range("E46").Value = "'  +1.23%  "
$ws.Range("D47").Value = "'9.330"
$ws.Range("E47").Value = "'  +2.16%  "
$ws.Range("D48").Value = "'7.110"
$ws.Range("E48").Value = "'  +2.67%  "
$ws.Range("D49").Value = "'35.67"
$ws.Range("E49").Value = "'  +4.76%  "
$ws.Range("D50").Value = "'921.76"
$ws.Range("E50").Value = "'  +0.31%  "
$ws.Range("D51").Value = "'0.3904"
$ws.Range("E51").Value = "'  +4.64%  "
